$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9300694554254023
$ws.Range("J2").Value = 0.9300694554254023
$ws.Range("M2").Value = 0.303146
$ws.Range("N2").Value = 0.909438
$ws.Range("O2").Value = 0.005142855213700541
$ws.Range("P2").Value = 0.005142855213700542
$ws.Range("Q2").Value = 0.7601753767146666
$ws.Range("R2").Value = 6.841578390432
$ws.Range("S2").Value = 0.004783212547938153
$ws.Range("T2").Value = 0.004783212547938153

# Row 3
$ws.Range("I3").Value = 0.9300694554254023
$ws.Range("J3").Value = 0.9300694554254023
$ws.Range("O3").Value = 0.2877784259203595
$ws.Range("P3").Value = 0.2877784259203595
$ws.Range("S3").Value = 0.2676539238789282
$ws.Range("T3").Value = 0.2676539238789282

# Row 4
$ws.Range("I4").Value = 0.9300694554254023
$ws.Range("J4").Value = 0.9300694554254023
$ws.Range("M4").Value = 41.67881
$ws.Range("N4").Value = 125.03643
$ws.Range("O4").Value = 0.7070787188659401
$ws.Range("P4").Value = 0.7070787188659401
$ws.Range("Q4").Value = 104.5146731039467
$ws.Range("R4").Value = 940.63205793552
$ws.Range("S4").Value = 0.657632318998536
$ws.Range("T4").Value = 0.657632318998536

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1885443333333333
$ws.Range("H5").Value = 0.5656329999999999
$ws.Range("I5").Value = 0.06993054457459773
$ws.Range("J5").Value = 0.06993054457459771
$ws.Range("M5").Value = 0.303146
$ws.Range("N5").Value = 0.909438
$ws.Range("O5").Value = 0.005142855213700541
$ws.Range("P5").Value = 0.005142855213700542
$ws.Range("Q5").Value = 0.05715646047266666
$ws.Range("R5").Value = 0.5144081442539999
$ws.Range("S5").Value = 0.000359642665762388
$ws.Range("T5").Value = 0.000359642665762388

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1885443333333333
$ws.Range("H6").Value = 0.5656329999999999
$ws.Range("I6").Value = 0.06993054457459773
$ws.Range("J6").Value = 0.06993054457459771
$ws.Range("O6").Value = 0.2877784259203595
$ws.Range("P6").Value = 0.2877784259203595
$ws.Range("Q6").Value = 3.198300465893889
$ws.Range("R6").Value = 28.78470419304499
$ws.Range("S6").Value = 0.02012450204143127
$ws.Range("T6").Value = 0.02012450204143126

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1885443333333333
$ws.Range("H7").Value = 0.5656329999999999
$ws.Range("I7").Value = 0.06993054457459773
$ws.Range("J7").Value = 0.06993054457459771
$ws.Range("M7").Value = 41.67881
$ws.Range("N7").Value = 125.03643
$ws.Range("O7").Value = 0.7070787188659401
$ws.Range("P7").Value = 0.7070787188659401
$ws.Range("Q7").Value = 7.858303445576666
$ws.Range("R7").Value = 70.72473101018998
$ws.Range("S7").Value = 0.04944639986740408
$ws.Range("T7").Value = 0.04944639986740407
